$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.834.27'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.636.90'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  -1.37%  '
$ws.Range("D5").Value = '''214.48'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '''0.5017'
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  -1.59%  '
$ws.Range("D8").Value = '''0.2559'
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("D9").Value = '''0.06346'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = '''19.35'
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").Value = '''0.07781'
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.695.60'
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.237'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.861.43'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '''0.5406'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '0.0₅7830'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").Value = '''64.13'
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = '25.874.46'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '''1.006'
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("D20").Value = '''194.97'
$ws.Range("E20").Value = '  -4.57%  '
$ws.Range("D21").Value = '''4.351'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").Value = '''9.842'
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("D23").Value = '''5.940'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '''1.005'
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("D25").Value = '''1.901'
$ws.Range("E25").Value = '  -4.08%  '
$ws.Range("D26").Value = '''139.77'
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").Value = '''0.1125'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '''6.775'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '''15.55'
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '''1.238'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '''0.04842'
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("D32").Value = '''3.232'
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").Value = '''3.163'
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").Value = '''1.526'
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").Value = '''2.373'
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '''0.8837'
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '''2.597'
$ws.Range("E37").Value = '  -1.47%  '
$ws.Range("D38").Value = '1.125.56'
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("D39").Value = '''0.5482'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("D40").Value = '''0.01557'
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = '''1.004'
$ws.Range("E41").Value = '  -1.70%  '
$ws.Range("D42").Value = '''5.659'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = '''0.8081'
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("D44").Value = '''99.30'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("E45").Value = '  +5.12%  '
$ws.Range("D46").Value = '1.773.02'
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").Value = '''0.4531'
$ws.Range("E47").Value = '  -0.88%  '
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("D49").Value = '''54.58'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("D50").Value = '''0.05042'
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  -1.17%  '
